# Fruta / hortaliza, semanal
# Updates columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) for rows 2-40 of the
# "Terminal La Palmera de La Serena - Coco" sheet. Row 41 is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: @(RowNumber, Fecha, Volumen, PrecioMinimo, PrecioMaximo, PrecioPromedioPonderado, PrecioPorKg)
$data = @(
    @(2, 44809, 60, 27000, 28000, 27500, 1375),
    @(3, 44333, 100, 19500, 20000, 19750, 988),
    @(4, 44880, 100, 28000, 30000, 29000, 1450),
    @(5, 44879, 100, 28000, 30000, 29000, 1450),
    @(6, 44810, 100, 27000, 28000, 27500, 1375),
    @(7, 44365, 100, 20000, 21000, 20500, 1025),
    @(8, 44784, 160, 27000, 28000, 27500, 1375),
    @(9, 44874, 240, 29000, 30000, 29500, 1475),
    @(10, 44428, 100, 20000, 21000, 20500, 1025),
    @(11, 44315, 100, 20000, 21000, 20500, 1025),
    @(12, 44335, 200, 19000, 20000, 19500, 975),
    @(13, 44445, 160, 20000, 21000, 20500, 1025),
    @(14, 44882, 120, 28000, 30000, 29000, 1450),
    @(15, 44301, 100, 18000, 19000, 18500, 925),
    @(16, 44473, 40, 19500, 20000, 19750, 988),
    @(17, 44778, 100, 23000, 24000, 23500, 1175),
    @(18, 44466, 100, 20000, 21000, 20500, 1025),
    @(19, 44427, 200, 20000, 21000, 20500, 1025),
    @(20, 44410, 200, 20000, 21000, 20500, 1025),
    @(21, 44467, 200, 20000, 21000, 20500, 1025),
    @(22, 44448, 100, 20000, 21000, 20500, 1025),
    @(23, 44350, 160, 19000, 20000, 19500, 975),
    @(24, 44434, 100, 20000, 21000, 20500, 1025),
    @(25, 44336, 100, 19500, 20000, 19750, 988),
    @(26, 44474, 200, 19000, 20000, 19500, 975),
    @(27, 44441, 160, 20000, 21000, 20500, 1025),
    @(28, 44782, 200, 23500, 24000, 23750, 1188),
    @(29, 44776, 160, 23000, 24000, 23500, 1175),
    @(30, 44418, 200, 20000, 21000, 20500, 1025),
    @(31, 44442, 140, 20000, 21000, 20500, 1025),
    @(32, 44462, 100, 19500, 20000, 19750, 988),
    @(33, 44326, 160, 19500, 20000, 19750, 988),
    @(34, 44435, 260, 20000, 22000, 21115, 1056),
    @(35, 44781, 160, 23000, 24000, 23500, 1175),
    @(36, 44407, 160, 20000, 21000, 20500, 1025),
    @(37, 44431, 160, 21000, 22000, 21500, 1075),
    @(38, 44420, 160, 20000, 21000, 20500, 1025),
    @(39, 44343, 100, 19500, 20000, 19750, 988),
    @(40, 44417, 160, 20000, 21000, 20500, 1025)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $row[2]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $row[3]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[4]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[5]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $row[6]   # S - Precio $/Kg
}
